$d = $word.ActiveDocument

# Edit 1: "anchor tag" paragraph - merge 3 runs (split by a gramStart/gramEnd
# proofErr pair) into a single run, and collapse the double space in
# "scroll  the" to a single space.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*anchor tag*") {
        $start = $p.Range.Start
        $end = $p.Range.End - 1
        $whole = $d.Range($start, $end)
        $whole.Delete()
        $ins = $d.Range($start, $start)
        $ins.InsertAfter("“#” is an anchor tag. It will suddenly scroll the person to a spot on the website.")
        break
    }
}

# Edit 2: "Wrapping an <img>" paragraph - merge 3 runs (split by a
# spellStart/spellEnd proofErr pair) into a single run.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Wrapping an*") {
        $start = $p.Range.Start
        $end = $p.Range.End - 1
        $whole = $d.Range($start, $end)
        $whole.Delete()
        $ins = $d.Range($start, $start)
        $ins.InsertAfter("Wrapping an <img> with <a></a> makes the logo or image clickable. ")
        break
    }
}

# Edit 3: "In-line CSS" paragraph - merge only the first 3 runs (split by a
# spellStart/spellEnd proofErr pair around "unreusable") into one run, while
# leaving the trailing "pages" run (wrapped in its own gramStart/gramEnd
# proofErr pair) untouched.
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*In-line CSS*") {
        $start = $p.Range.Start
        $idx = $t.IndexOf("pages")
        $boundary = $start + $idx
        # Insert the merged replacement text right before "pages" first, then
        # delete the original (still-split) prefix. Doing it in this order
        # keeps the new run outside/before the proofErr markers that wrap
        # "pages" instead of landing inside them.
        $insPoint = $d.Range($boundary, $boundary)
        $insPoint.InsertBefore("In-line CSS: not favored because it is unreusable on multiple ")
        $oldPrefix = $d.Range($start, $boundary)
        $oldPrefix.Delete()
        break
    }
}
